$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert new row 28 with the new entry, pushing existing row 27 data down is NOT needed
# since row 28 in the sheet was previously empty (rows jump from 27 to 34).
$ws.Cells.Item(28, 1).Value = (Get-Date -Year 2013 -Month 6 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(28, 2).Value = "Projekthandbuch, Statusbericht, Analyse workbench.xmi"
$ws.Cells.Item(28, 3).Value = 8

# Apply same styles as row 27 (date style s=3 for A, s=1 for B and C)
$ws.Cells.Item(27, 1).Copy()
$ws.Cells.Item(28, 1).PasteSpecial(-4122) # xlPasteFormats
$ws.Cells.Item(27, 2).Copy()
$ws.Cells.Item(28, 2).PasteSpecial(-4122)
$ws.Cells.Item(27, 3).Copy()
$ws.Cells.Item(28, 3).PasteSpecial(-4122)

# Update selection to B28
$ws.Range("B28").Select()

$wb.Save()
